$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.098.92'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').Value = '3.278.42'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'575.76"
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = "'179.23"
$ws.Range('E6').Value = '  -3.59%  '
$ws.Range('D7').Value = "'0.624"
$ws.Range('E7').Value = '  +2.83%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -3.01%  '
$ws.Range('D10').Value = "'6.72"
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').Value = '3.845.47'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D14').Value = '66.161.83'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').Value = "'0.0000163"
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('D17').Value = '3.272.42'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('D18').Value = "'432.74"
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('E21').Value = '  -4.10%  '
$ws.Range('D22').Value = "'72.00"
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '3.412.35'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').Value = "'0.506"
$ws.Range('E25').Value = '  -1.57%  '
$ws.Range('E26').Value = '  +3.31%  '
$ws.Range('E27').Value = '  -5.16%  '
$ws.Range('D28').Value = "'8.89"
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('D31').Value = "'22.28"
$ws.Range('E31').Value = '  -2.86%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = "'5.16"
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('D34').Value = "'6.59"
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('E35').Value = '  -4.49%  '
$ws.Range('D36').Value = "'157.09"
$ws.Range('E36').Value = '  -3.19%  '
$ws.Range('E37').Value = '  -5.77%  '
$ws.Range('D38').Value = "'26.61"
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('D40').Value = '2.771.54'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('D41').Value = "'0.776"
$ws.Range('E41').Value = '  -1.84%  '
$ws.Range('E42').Value = '  -3.86%  '
$ws.Range('D43').Value = "'40.27"
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = "'6.02"
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('E45').Value = '  -2.29%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'321.63"
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = "'2.30"
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('D48').Value = "'23.29"
$ws.Range('E48').Value = '  -5.92%  '
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('D51').Value = "'0.999"
$ws.Range('E51').Value = '  -0.03%  '
